$wb = $excel.ActiveWorkbook

# Updated "rit" (item-rest correlation) values for column C, rows 2-5, per worksheet.
# Rerun of distractor analysis for simulated example (use_wle = TRUE).
$newValues = @{}
$newValues["reg70001"] = @(-0.24,-0.266,0.467,-0.246)
$newValues["reg70003"] = @(-0.31,0.548,-0.199,-0.222)
$newValues["reg70005"] = @(-0.167,-0.122,-0.092,0.435)
$newValues["reg70006"] = @(-0.251,0.458,-0.253,-0.221)
$newValues["reg70007"] = @(0.45,-0.207,-0.256,-0.233)
$newValues["reg70008"] = @(0.531,-0.318,-0.239,-0.23)
$newValues["reg70009"] = @(-0.236,0.51,-0.26,-0.252)
$newValues["reg70011"] = @(0.537,-0.279,-0.251,-0.262)
$newValues["reg70012"] = @(-0.248,-0.268,0.544,-0.244)
$newValues["reg70013"] = @(-0.243,-0.211,-0.277,0.54)
$newValues["reg70014"] = @(-0.24,-0.177,0.518,-0.262)
$newValues["reg70015"] = @(0.511,-0.198,-0.204,-0.241)
$newValues["reg70017"] = @(-0.177,-0.208,0.537,-0.268)
$newValues["reg70018"] = @(-0.185,-0.183,-0.243,0.514)
$newValues["reg70020"] = @(-0.133,-0.196,-0.248,0.56)
$newValues["reg70021"] = @(0.521,-0.171,-0.172,-0.159)

foreach ($name in $newValues.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $vals = $newValues[$name]
    for ($r = 2; $r -le 5; $r++) {
        $ws.Cells.Item($r, 3).Value = $vals[$r - 2]
    }
}
